# The Akan sentence "Oforo, foro, foro a, twon ! Na ode aka no bio, nso
# oforo, foro, foro, twon!" (row 6, column A) and its English translation
# "He climbed, and climbed, and climbed; in vain. He strove again, again
# he made to climb, and climb, and climb; in vain." (row 6, column B) were
# originally stored as single combined sentences. This change splits each
# of them into two separate rows so that each clause gets its own row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row below the row that currently holds the combined
# sentence (row 6), so that we end up with two rows for the two clauses.
$ws.Rows("7").Insert()

# First clause of the split sentence stays on row 6, second clause goes
# on the newly inserted row 7. Fill in column A (Akan) for both rows
# first, then column B (English) for both rows, to match the order in
# which the new shared strings were appended to the workbook.
$ws.Range("A6").Value2 = "Oforo, foro, foro a, twon !"
$ws.Range("A7").Value2 = "Na ode aka no bio, nso oforo, foro, foro, twon!"
$ws.Range("B6").Value2 = "He climbed, and climbed, and climbed; in vain."
$ws.Range("B7").Value2 = "He strove again, again he made to climb, and climb, and climb; in vain."

# Let Excel recompute the row height for row 6 now that it holds shorter
# text (row 7 already has no explicit height since it is brand new).
$ws.Rows("6").AutoFit()

# Update the selection to match the saved state of the workbook.
[void]$ws.Range("E7").Select()
